# Formed the consolidated report
# Recalculate the "Absent" column (H) as the complement of the "Real" column (E):
#   Absent = 1 - Real
# This fixes rows where H was stale/blank so it is consistent with column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 3) { $lastRow = 21 }

for ($r = 3; $r -le $lastRow; $r++) {
    $realCell = $ws.Cells.Item($r, 5)   # column E - Real
    $absentCell = $ws.Cells.Item($r, 8) # column H - Absent

    $realValue = $realCell.Value2
    if ($null -eq $realValue) { $realValue = 0 }

    $absentCell.Value2 = 1 - $realValue
}
